$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 131.24
$ws.Range("K5").Value = 69
$ws.Range("L5").Value = 170.07
